# Update cryptos list snapshot (price/volume refresh + Dai/BitcoinCash row swap)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Bitcoin)
$ws.Range("D2").Value = "26.671.76"
$ws.Range("E2").Value = "  -0.06%  "
# Row 3 (Ethereum)
$ws.Range("D3").Value = "1.597.26"
$ws.Range("E3").Value = "  -0.21%  "
# Row 4 (TetherUSD)
$ws.Range("E4").Value = "  +0.34%  "
# Row 5 (BNB)
$r = $ws.Range("D5")
$origStyle = $r.Style
$r.NumberFormat = "@"  # preserve text type for numeric-looking value
$r.Value = "211.56"
$r.Style = $origStyle
$ws.Range("E5").Value = "  +0.31%  "
# Row 6 (XRP)
$ws.Range("E6").Value = "  -0.07%  "
# Row 7 (USDC)
$ws.Range("E7").Value = "  +0.35%  "
# Row 8 (Dogecoin)
$ws.Range("E8").Value = "  +0.05%  "
# Row 9 (Cardano)
$r = $ws.Range("D9")
$origStyle = $r.Style
$r.NumberFormat = "@"  # preserve text type for numeric-looking value
$r.Value = "0.247"
$r.Style = $origStyle
$ws.Range("E9").Value = "  +0.45%  "
# Row 10 (Solana)
$r = $ws.Range("D10")
$origStyle = $r.Style
$r.NumberFormat = "@"  # preserve text type for numeric-looking value
$r.Value = "19.48"
$r.Style = $origStyle
$ws.Range("E10").Value = "  -0.64%  "
# Row 11 (TRON)
$ws.Range("E11").Value = "  -0.07%  "
# Row 12 (WrappedliquidstakedEther2.0)
$ws.Range("D12").Value = "1.821.51"
$ws.Range("E12").Value = "  -0.04%  "
# Row 13 (WrappedEther)
$ws.Range("D13").Value = "1.583.69"
$ws.Range("E13").Value = "  -0.68%  "
# Row 14 (Polkadot)
$r = $ws.Range("D14")
$origStyle = $r.Style
$r.NumberFormat = "@"  # preserve text type for numeric-looking value
$r.Value = "4.03"
$r.Style = $origStyle
$ws.Range("E14").Value = "  +0.06%  "
# Row 15 (Polygon)
$ws.Range("E15").Value = "  +0.44%  "
# Row 16 (Litecoin)
$r = $ws.Range("D16")
$origStyle = $r.Style
$r.NumberFormat = "@"  # preserve text type for numeric-looking value
$r.Value = "65.02"
$r.Style = $origStyle
$ws.Range("E16").Value = "  +0.35%  "
# Row 17 (WrappedBTC)
$ws.Range("D17").Value = "26.645.53"
$ws.Range("E17").Value = "  -0.09%  "
# Row 18 (ShibaInu)
$ws.Range("D18").Value = "0.0₃0739"
$ws.Range("E18").Value = "  +1.43%  "
# Row 19 (Dai)
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$r = $ws.Range("D19")
$origStyle = $r.Style
$r.NumberFormat = "@"  # preserve text type for numeric-looking value
$r.Value = "209.25"
$r.Style = $origStyle
$ws.Range("E19").Value = "  +0.14%  "
# Row 20 (BitcoinCash)
$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$r = $ws.Range("D20")
$origStyle = $r.Style
$r.NumberFormat = "@"  # preserve text type for numeric-looking value
$r.Value = "1.00"
$r.Style = $origStyle
$ws.Range("E20").Value = "  +0.25%  "
# Row 21 (Chainlink)
$ws.Range("E21").Value = "  +4.28%  "
# Row 22 (Uniswap)
$ws.Range("E22").Value = "  +0.84%  "
# Row 23 (Toncoin)
$ws.Range("E23").Value = "  +1.63%  "
# Row 24 (Avalanche)
$r = $ws.Range("D24")
$origStyle = $r.Style
$r.NumberFormat = "@"  # preserve text type for numeric-looking value
$r.Value = "8.99"
$r.Style = $origStyle
$ws.Range("E24").Value = "  +0.90%  "
# Row 25 (Monero)
$r = $ws.Range("D25")
$origStyle = $r.Style
$r.NumberFormat = "@"  # preserve text type for numeric-looking value
$r.Value = "143.20"
$r.Style = $origStyle
$ws.Range("E25").Value = "  -1.64%  "
# Row 26 (BinanceUSD)
$ws.Range("E26").Value = "  +0.30%  "
# Row 27 (Cosmos)
$r = $ws.Range("D27")
$origStyle = $r.Style
$r.NumberFormat = "@"  # preserve text type for numeric-looking value
$r.Value = "7.11"
$r.Style = $origStyle
$ws.Range("E27").Value = "  -1.72%  "
# Row 28 (Stellar)
$ws.Range("E28").Value = "  -1.11%  "
# Row 29 (EthereumClassic)
$ws.Range("E29").Value = "  +0.10%  "
# Row 30 (Hedera)
$ws.Range("E30").Value = "  +1.74%  "
# Row 31 (PancakeSwap)
$ws.Range("E31").Value = "  +0.03%  "
# Row 32 (Filecoin)
$ws.Range("E32").Value = "  -0.06%  "
# Row 33 (InternetComputer(DFINITY))
$ws.Range("E33").Value = "  +0.41%  "
# Row 34 (Maker)
$ws.Range("D34").Value = "1.291.00"
$ws.Range("E34").Value = "  -0.18%  "
# Row 35 (ImmutableX)
$r = $ws.Range("D35")
$origStyle = $r.Style
$r.NumberFormat = "@"  # preserve text type for numeric-looking value
$r.Value = "0.618"
$r.Style = $origStyle
$ws.Range("E35").Value = "  -5.31%  "
# Row 36 (HuobiToken)
$ws.Range("E36").Value = "  +0.36%  "
# Row 37 (LidoDAOToken)
$ws.Range("E37").Value = "  -0.31%  "
# Row 38 (VeChain)
$r = $ws.Range("D38")
$origStyle = $r.Style
$r.NumberFormat = "@"  # preserve text type for numeric-looking value
$r.Value = "0.0171"
$r.Style = $origStyle
$ws.Range("E38").Value = "  -0.49%  "
# Row 39 (ARBITRUM)
$r = $ws.Range("D39")
$origStyle = $r.Style
$r.NumberFormat = "@"  # preserve text type for numeric-looking value
$r.Value = "0.827"
$r.Style = $origStyle
$ws.Range("E39").Value = "  -2.34%  "
# Row 40 (WEMIXToken)
$r = $ws.Range("D40")
$origStyle = $r.Style
$r.NumberFormat = "@"  # preserve text type for numeric-looking value
$r.Value = "1.05"
$r.Style = $origStyle
$ws.Range("E40").Value = "  +16.74%  "
# Row 41 (FraxShare)
$ws.Range("E41").Value = "  +0.82%  "
# Row 42 (MXToken)
$ws.Range("E42").Value = "  -0.56%  "
# Row 43 (TrustWalletToken)
$ws.Range("E43").Value = "  -0.78%  "
# Row 44 (Aave)
$r = $ws.Range("D44")
$origStyle = $r.Style
$r.NumberFormat = "@"  # preserve text type for numeric-looking value
$r.Value = "63.18"
$r.Style = $origStyle
$ws.Range("E44").Value = "  -1.01%  "
# Row 45 (RocketPoolETH)
$ws.Range("D45").Value = "1.733.15"
$ws.Range("E45").Value = "  -0.12%  "
# Row 46 (Quant)
$r = $ws.Range("D46")
$origStyle = $r.Style
$r.NumberFormat = "@"  # preserve text type for numeric-looking value
$r.Value = "90.99"
$r.Style = $origStyle
$ws.Range("E46").Value = "  +0.98%  "
# Row 47 (RenderToken)
$ws.Range("E47").Value = "  -2.73%  "
# Row 48 (Algorand)
$ws.Range("E48").Value = "  +1.14%  "
# Row 49 (Cronos)
$ws.Range("E49").Value = "  +0.91%  "
# Row 50 (USDD)
$ws.Range("E50").Value = "  +0.45%  "
# Row 51 (EnergySwap)
$r = $ws.Range("D51")
$origStyle = $r.Style
$r.NumberFormat = "@"  # preserve text type for numeric-looking value
$r.Value = "7.35"
$r.Style = $origStyle
$ws.Range("E51").Value = "  -1.39%  "
